# Generate Report for Handoff
# Update the "Latest Handback DateTime" column (D) for the most recently
# handed-back file (3338663b-41cd-4af3-8cd2-89dee62ff182) on both the
# zh-cn and de-de localization status sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value2 = "2016-03-09 12:39:55"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value2 = "2016-03-09 12:40:05"
